# Update master to output generated at 503736d
$d = $word.ActiveDocument

# --- Title date line ---------------------------------------------------
$d.Content.Find.Execute("2025-02-21 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-22 Saturday", 2) | Out-Null

# --- Row 1: replace the whole set of problems ---------------------------
# Original row 1 cells: 65÷9=7, 2 | 59÷8=7, 3 | 96÷7=13, 5 | 42÷7=6, 0 | 20÷8=2, 4
# Target row 1 cells:   42÷7=6, 0 | 75÷8=9, 3 | 10÷2=5, 0 | 62÷2=31, 0 | 39÷3=13, 0
# The first three cells are dropped, the fourth survives unchanged, the
# fifth is re-worded, and three brand-new cells are appended -- net
# effect: the row keeps 5 cells but with entirely different content.
# Cell-level insert/delete isn't available on this host, so the row is
# rebuilt instead: insert a fresh row before the old one (this inherits
# the old row's cell/run formatting), populate its five cells with the
# final values, then delete the stale row.
$t = $d.Tables.Item(1)
$oldRow1 = $t.Rows.Item(1)
$newRow1 = $t.Rows.Add($oldRow1)
$newRow1.Cells.Item(1).Range.Text = "42÷7=6, 0"
$newRow1.Cells.Item(2).Range.Text = "75÷8=9, 3"
$newRow1.Cells.Item(3).Range.Text = "10÷2=5, 0"
$newRow1.Cells.Item(4).Range.Text = "62÷2=31, 0"
$newRow1.Cells.Item(5).Range.Text = "39÷3=13, 0"
$t.Rows.Item(2).Delete()

# --- Row 5: simple in-place text swaps ----------------------------------
$d.Content.Find.Execute("69÷6=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "54÷5=10, 4", 2) | Out-Null
$d.Content.Find.Execute("70÷5=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷6=10, 2", 2) | Out-Null
$d.Content.Find.Execute("31÷3=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=3, 5", 2) | Out-Null
$d.Content.Find.Execute("94÷3=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "77÷9=8, 5", 2) | Out-Null
$d.Content.Find.Execute("64÷4=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "25÷3=8, 1", 2) | Out-Null

# --- Row 9: simple in-place text swaps -----------------------------------
$d.Content.Find.Execute("74÷6=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "82÷2=41, 0", 2) | Out-Null
$d.Content.Find.Execute("15÷8=1, 7", $true, $false, $false, $false, $false, $true, 1, $false, "57÷2=28, 1", 2) | Out-Null
$d.Content.Find.Execute("66÷4=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "69÷7=9, 6", 2) | Out-Null
$d.Content.Find.Execute("54÷3=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "39÷9=4, 3", 2) | Out-Null
$d.Content.Find.Execute("89÷6=14, 5", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=5, 7", 2) | Out-Null

# --- Row 13: simple in-place text swaps ----------------------------------
$d.Content.Find.Execute("36÷5=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "30÷3=10, 0", 2) | Out-Null
$d.Content.Find.Execute("48÷6=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "58÷2=29, 0", 2) | Out-Null
$d.Content.Find.Execute("17÷4=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "91÷9=10, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷8=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "92÷6=15, 2", 2) | Out-Null
$d.Content.Find.Execute("42÷9=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "18÷4=4, 2", 2) | Out-Null

# --- Row 17: simple in-place text swaps ----------------------------------
$d.Content.Find.Execute("41÷8=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "25÷4=6, 1", 2) | Out-Null
$d.Content.Find.Execute("97÷4=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "33÷4=8, 1", 2) | Out-Null
$d.Content.Find.Execute("89÷2=44, 1", $true, $false, $false, $false, $false, $true, 1, $false, "23÷9=2, 5", 2) | Out-Null
$d.Content.Find.Execute("92÷4=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "60÷9=6, 6", 2) | Out-Null
$d.Content.Find.Execute("78÷9=8, 6", $true, $false, $false, $false, $false, $true, 1, $false, "39÷5=7, 4", 2) | Out-Null
